$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "C:Antimalarial Use" row (original row 4) and the "I:Anemia" row (original row 7).
# Delete from bottom up so row indices don't shift unexpectedly.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(4).Delete()

# Rename two variable names that changed.
$ws.Range("B10").Value = "hlthst_duration_fctb_clst"
$ws.Range("B11").Value = "wtrdist_fctb_clst"

# Update the selected cell to match the saved view state.
$ws.Range("B11").Select()
